$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "HO" (sheet2.xml): insert a "Month" column (B) and a "Diiference"
# column (G) that subtracts this year's amount from last year's amount, then
# add a block of blank (but formatted) rows below the single data row.
# ---------------------------------------------------------------------------
$wsHO = $wb.Worksheets.Item("HO")

# Insert new column B ("Month") - shifts old B:E to C:F
$wsHO.Columns.Item(2).Insert()
# Insert new column G ("Diiference") after the (now shifted) data columns
$wsHO.Columns.Item(7).Insert()

$wsHO.Range("B1").Value = "Month"
$wsHO.Range("B2").Value = "December"
$wsHO.Range("G1").Value = "Diiference"
$wsHO.Range("G2").Formula = "=F2-D2"

$wsHO.Columns.Item(2).ColumnWidth = 46.25
$wsHO.Columns.Item(7).ColumnWidth = 8.92

# Add the trailing blank rows (3-22) so column B carries the same formatting
# all the way down as in the finished report.
$wsHO.Range("B3:B22").Style = $wsHO.Range("A2").Style

# ---------------------------------------------------------------------------
# Sheet "New Stores" (sheet3.xml): insert a "Month" column (B) populated with
# "December" for every data row.
# ---------------------------------------------------------------------------
$wsNew = $wb.Worksheets.Item("New Stores")
$wsNew.Columns.Item(2).Insert()
$wsNew.Range("B1").Value = "Month"
$wsNew.Range("B2").Value = "December"
$wsNew.Range("B3").Value = "December"
$wsNew.Range("B4").Value = "December"
$wsNew.Range("B5").Value = "December"
$wsNew.Range("B6").Value = "December"
$wsNew.Columns.Item(2).ColumnWidth = 41.59

# ---------------------------------------------------------------------------
# Sheet "Closed Stores" (sheet4.xml): insert a "Month" column (B) populated
# with "December" for every data row.
# ---------------------------------------------------------------------------
$wsClosed = $wb.Worksheets.Item("Closed Stores")
$wsClosed.Columns.Item(2).Insert()
$wsClosed.Range("B1").Value = "Month"
$wsClosed.Range("B2").Value = "December"
$wsClosed.Range("B3").Value = "December"
$wsClosed.Range("B4").Value = "December"
$wsClosed.Columns.Item(2).ColumnWidth = 29.25

# ---------------------------------------------------------------------------
# Selections / scroll position per sheet. Do the sheet that should stay the
# active tab (YOY) LAST, so it ends up as the selected tab on save.
# ---------------------------------------------------------------------------
$wsHO.Range("G1:G2").Select()
$wsNew.Range("D1").Select()
$wsClosed.Range("E1").Select()

$wsYOY = $wb.Worksheets.Item(1)
$wsYOY.Range("B15").Select()
